$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Replace the old (now-dead) share-me-sanity login URL with the
#    new merny.netlify.app/auth URL.
# ------------------------------------------------------------------
$d.Content.Find.Execute("https://share-me-sanity.netlify.app/login", $true, $false, $false, $false, $false,
                         $true, 1, $false, "https://merny.netlify.app/auth", 2) | Out-Null

# ------------------------------------------------------------------
# 2) Word stamps a "_GoBack" bookmark (empty, zero-width) right where
#    the last edit happened -- here, immediately after the edited
#    run, at the end of the paragraph's text (before the paragraph
#    mark). Rebuild the paragraph's WordOpenXML with the bookmark
#    appended after the run to reproduce that marker precisely.
# ------------------------------------------------------------------
$para = $d.Paragraphs(1)
$rng = $para.Range

$bookmarkPkg = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>https://merny.netlify.app/auth</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rng.InsertXML($bookmarkPkg) | Out-Null

# InsertXML above swaps the paragraph's run content for a fresh <w:p>,
# which leaves the original (now empty) trailing paragraph mark
# behind as a spurious extra paragraph. Merge that stray paragraph
# mark back into the (only) paragraph.
if ($d.Paragraphs.Count -gt 1) {
    $p1 = $d.Paragraphs(1)
    $tailRng = $d.Range($p1.Range.End - 1, $d.Content.End)
    $tailRng.Delete() | Out-Null
}

# ------------------------------------------------------------------
# 3) Section/page setup: page size Letter (was A4), and header /
#    footer distance + column spacing of 36pt (720 twips, was 708).
# ------------------------------------------------------------------
$ps = $d.PageSetup
$ps.PageWidth = 612
$ps.PageHeight = 792
$ps.HeaderDistance = 36
$ps.FooterDistance = 36
$ps.TextColumns.Spacing = 36
